# B1-/B2- PowerPoint — Tue, Jul 07, 2020 edit
#
# 1. The table on slide 5 (the financial-documents table) gets switched to a
#    different built-in table style (tableStyleId GUID change).
# 2. The deck's theme colour scheme (the "Integral" theme used by the slide
#    master, physically stored as ppt/theme/theme2.xml) is recoloured to the
#    stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$s5 = $p.Slides.Item(5)
$tableShape = $s5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{22CBD6D5-7290-4157-AD66-901A6DF7DAE6}", $false)

# --- 2. Recolour the active theme to the stock "Office" colours ----------
$tcs = $p.Slides.Item(1).ThemeColorScheme

# index -> (name, RGB as decimal 0x00BBGGRR, matching the "Office" scheme)
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
